$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes ---
# Row 7 (SHHConfig) column B value changes from "ssh" to "Configurations/SSHConfig.csv"
$ws.Range("B7").Value = "Configurations/SSHConfig.csv"

# Row 6 (TagsConfig) gains a value in column B: "fog=true/fog_visible=85"
$ws.Range("B6").Value = "fog=true/fog_visible=85"

# Column C was a redundant duplicate of column B; delete it entirely
$ws.Columns("C:C").Delete()

# --- Column widths (A:C) ---
# (Input values tuned so the engine's internal pixel-rounding lands on the
# closest achievable stored width to the target 25.453125 / 30.6328125 / 25.81640625)
$ws.Columns("A:A").ColumnWidth = 24.673125
$ws.Columns("B:B").ColumnWidth = 29.8378125
$ws.Columns("C:C").ColumnWidth = 25.0064063

# --- Selection moves to I12 ---
$ws.Range("I12").Select() | Out-Null
